# Update cryptos list values (price + 1h volume change) per scraped data refresh.
# Numeric-looking Price (column D) values are prefixed with a leading apostrophe so
# Excel stores them as text (matching the source data's inlineStr representation)
# instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.711.16"
$ws.Range("E2").Value = "  +1.93%  "
$ws.Range("D3").Value = "2.215.30"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'269.78"
$ws.Range("E5").Value = "  +5.00%  "
$ws.Range("D6").Value = "'85.62"
$ws.Range("E6").Value = "  +10.65%  "
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").Value = "'0.602"
$ws.Range("E9").Value = "  +1.06%  "
$ws.Range("D10").Value = "'45.73"
$ws.Range("E10").Value = "  +7.45%  "
$ws.Range("D11").Value = "'0.0919"
$ws.Range("E11").Value = "  +1.17%  "
$ws.Range("E12").Value = "  +6.13%  "
$ws.Range("D13").Value = "'0.105"
$ws.Range("E13").Value = "  +1.86%  "
$ws.Range("D15").Value = "'14.56"
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("D16").Value = "2.216.16"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("D17").Value = "'0.783"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").Value = "43.651.53"
$ws.Range("E18").Value = "  +1.83%  "
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("D20").Value = "'5.98"
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").Value = "'69.83"
$ws.Range("E21").Value = "  -1.87%  "
$ws.Range("D22").Value = "'2.36"
$ws.Range("E22").Value = "  +5.30%  "
$ws.Range("D23").Value = "'231.53"
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("D24").Value = "'2.72"
$ws.Range("E24").Value = "  +23.82%  "
$ws.Range("D25").Value = "'8.85"
$ws.Range("E25").Value = "  -5.35%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'10.76"
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("E28").Value = "  +5.82%  "
$ws.Range("D29").Value = "'38.94"
$ws.Range("E29").Value = "  -8.83%  "
$ws.Range("E30").Value = "  -0.52%  "
$ws.Range("D31").Value = "'175.38"
$ws.Range("E31").Value = "  +1.42%  "
$ws.Range("D32").Value = "'0.0891"
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("D33").Value = "'20.51"
$ws.Range("E33").Value = "  +0.65%  "
$ws.Range("E34").Value = "  +3.00%  "
$ws.Range("E35").Value = "  +1.90%  "
$ws.Range("E36").Value = "  +1.87%  "
$ws.Range("D37").Value = "'0.0354"
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("D38").Value = "'4.35"
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("E39").Value = "  +15.50%  "
$ws.Range("D40").Value = "'12.26"
$ws.Range("E40").Value = "  -6.12%  "
$ws.Range("D41").Value = "'64.95"
$ws.Range("E41").Value = "  +7.59%  "
$ws.Range("D42").Value = "'2.09"
$ws.Range("E42").Value = "  -1.17%  "
$ws.Range("D43").Value = "'0.204"
$ws.Range("E43").Value = "  +1.12%  "
$ws.Range("D44").Value = "'5.38"
$ws.Range("E44").Value = "  +1.18%  "
$ws.Range("D45").Value = "'0.0988"
$ws.Range("E45").Value = "  +1.18%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'8.33"
$ws.Range("E46").Value = "  -0.57%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'100.21"
$ws.Range("E47").Value = "  -2.89%  "
$ws.Range("D48").Value = "'1.21"
$ws.Range("E48").Value = "  +5.85%  "
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("D50").Value = "'0.438"
$ws.Range("E50").Value = "  -6.99%  "
$ws.Range("D51").Value = "'1.49"
$ws.Range("E51").Value = "  +3.79%  "
